$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 58926.117
$ws.Range("I33").Value = 66769.60000000001
$ws.Range("K33").Value = 66769.60000000001
$ws.Range("M33").Value = -66540.60000000001

# Row 48
$ws.Range("H48").Value = 8223.799999999999
$ws.Range("J48").Value = 8223.799999999999
$ws.Range("L48").Value = 24671.4
$ws.Range("N48").Value = -25255.4

# Row 51
$ws.Range("H51").Value = 2847.6667
$ws.Range("I51").Value = 3520
$ws.Range("J51").Value = 2175.3333
$ws.Range("K51").Value = 3520
$ws.Range("L51").Value = 2175.3333
$ws.Range("M51").Value = -3036
$ws.Range("N51").Value = -3143.3333

# Row 56
$ws.Range("H56").Value = 8223.799999999999
$ws.Range("J56").Value = 8223.799999999999
$ws.Range("L56").Value = 24671.4
$ws.Range("N56").Value = -25739.4

# Row 74
$ws.Range("H74").Value = 6492490.5
$ws.Range("I74").Value = 6492490.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 6492490.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -6491554.5
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 6492490.5
$ws.Range("I77").Value = 6492490.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 32462452.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -32457772.5
$ws.Range("N77").ClearContents()

# Row 138
$ws.Range("H138").Value = 2181.63
$ws.Range("I138").Value = 1414.7778
$ws.Range("J138").Value = 2612.9844
$ws.Range("K138").Value = 4244.3334
$ws.Range("L138").Value = 7838.9532
$ws.Range("M138").Value = 895.6665999999996
$ws.Range("N138").Value = -18118.9532


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 103606.8
$ws.Range("I2").Value = 168969.67
$ws.Range("K2").Value = 168969.67
$ws.Range("M2").Value = -168856.67

# Row 45
$ws.Range("H45").Value = 2781.0667
$ws.Range("I45").Value = 1610.5454
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 1610.5454
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -1233.5454
$ws.Range("N45").Value = -6754

# Row 116
$ws.Range("H116").Value = 103606.8
$ws.Range("I116").Value = 168969.67
$ws.Range("K116").Value = 168969.67
$ws.Range("M116").Value = -166675.67


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 103606.8
$ws.Range("I3").Value = 168969.67
$ws.Range("K3").Value = 168969.67
$ws.Range("M3").Value = -168855.67

# Row 105
$ws.Range("H105").Value = 3558.9473
$ws.Range("I105").Value = 2090
$ws.Range("K105").Value = 2090
$ws.Range("M105").Value = -343

# Row 134
$ws.Range("H134").Value = 35729.39
$ws.Range("I134").Value = 1871.2273
$ws.Range("J134").Value = 135046.67
$ws.Range("K134").Value = 5613.6819
$ws.Range("L134").Value = 405140.01
$ws.Range("M134").Value = -3078.6819
$ws.Range("N134").Value = -410210.01


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 835.55
$ws.Range("I68").Value = 516.9677
$ws.Range("J68").Value = 1355.3422
$ws.Range("K68").Value = 1550.9031
$ws.Range("L68").Value = 4066.0266
$ws.Range("M68").Value = -739.9031
$ws.Range("N68").Value = -5688.0266

# Row 71
$ws.Range("H71").Value = 835.55
$ws.Range("I71").Value = 516.9677
$ws.Range("J71").Value = 1355.3422
$ws.Range("K71").Value = 4652.7093
$ws.Range("L71").Value = 12198.0798
$ws.Range("M71").Value = -596.7093000000004
$ws.Range("N71").Value = -20310.0798

# Row 129
$ws.Range("H129").Value = 30510.115
$ws.Range("J129").Value = 35407.465
$ws.Range("L129").Value = 106222.395
$ws.Range("N129").Value = -116222.395


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 35997
$ws.Range("I12").Value = 1980
$ws.Range("J12").Value = 47336
$ws.Range("K12").Value = 1980
$ws.Range("L12").Value = 47336
$ws.Range("M12").Value = -1840
$ws.Range("N12").Value = -47616

# Row 138
$ws.Range("H138").Value = 39831.11
$ws.Range("J138").Value = 39831.11
$ws.Range("L138").Value = 39831.11
$ws.Range("N138").Value = -50111.11


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1155.875
$ws.Range("I7").Value = 1083.1666
$ws.Range("J7").Value = 1199.5
$ws.Range("K7").Value = 1083.1666
$ws.Range("L7").Value = 1199.5
$ws.Range("M7").Value = -971.1666
$ws.Range("N7").Value = -1423.5

# Row 36
$ws.Range("H36").Value = 38333.332
$ws.Range("J36").Value = 38333.332
$ws.Range("L36").Value = 38333.332
$ws.Range("N36").Value = -39457.332

# Row 40
$ws.Range("H40").Value = 50002
$ws.Range("I40").Value = 50002
$ws.Range("K40").Value = 50002
$ws.Range("M40").Value = -49866

# Row 46
$ws.Range("H46").Value = 499.375
$ws.Range("I46").Value = 372.85715
$ws.Range("J46").Value = 597.7778
$ws.Range("K46").Value = 372.85715
$ws.Range("L46").Value = 597.7778
$ws.Range("M46").Value = -184.85715
$ws.Range("N46").Value = -973.7778

# Row 61
$ws.Range("H61").Value = 1112.6875
$ws.Range("I61").Value = 930.4
$ws.Range("J61").Value = 1416.5
$ws.Range("K61").Value = 930.4
$ws.Range("L61").Value = 1416.5
$ws.Range("M61").Value = -728.4
$ws.Range("N61").Value = -1820.5

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 68
$ws.Range("H68").Value = 1877.7778
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 1725
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 1725
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -3223

# Row 71
$ws.Range("H71").Value = 1877.7778
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 1725
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 8625
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -16113

# Row 113
$ws.Range("H113").Value = 1112.6875
$ws.Range("I113").Value = 930.4
$ws.Range("J113").Value = 1416.5
$ws.Range("K113").Value = 930.4
$ws.Range("L113").Value = 1416.5
$ws.Range("M113").Value = 1239.6
$ws.Range("N113").Value = -5756.5

# Row 126
$ws.Range("H126").Value = 1155.875
$ws.Range("I126").Value = 1083.1666
$ws.Range("J126").Value = 1199.5
$ws.Range("K126").Value = 3249.4998
$ws.Range("L126").Value = 3598.5
$ws.Range("M126").Value = -779.4998000000001
$ws.Range("N126").Value = -8538.5

# Row 136
$ws.Range("H136").Value = 2121.4238
$ws.Range("I136").Value = 1382.3611
$ws.Range("K136").Value = 4147.0833
$ws.Range("M136").Value = -1597.0833

